$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H18").Value = 360.125
$ws_ALC.Range("I18").Value = 360.125
$ws_ALC.Range("K18").Value = 360.125
$ws_ALC.Range("M18").Value = -76.125
$ws_ALC.Range("H38").Value = 1415.0625
$ws_ALC.Range("I38").Value = 1220.0834
$ws_ALC.Range("J38").Value = 2000
$ws_ALC.Range("K38").Value = 3660.2502
$ws_ALC.Range("L38").Value = 6000
$ws_ALC.Range("M38").Value = -3288.2502
$ws_ALC.Range("N38").Value = -6744
$ws_ALC.Range("H43").Value = 632.3570999999999
$ws_ALC.Range("I43").Value = 767.2857
$ws_ALC.Range("K43").Value = 767.2857
$ws_ALC.Range("M43").Value = -698.2857
$ws_ALC.Range("H74").Value = 3695.4119
$ws_ALC.Range("I74").Value = 3617.8
$ws_ALC.Range("J74").Value = 3806.2856
$ws_ALC.Range("K74").Value = 3617.8
$ws_ALC.Range("L74").Value = 3806.2856
$ws_ALC.Range("M74").Value = -2681.8
$ws_ALC.Range("N74").Value = -5678.2856
$ws_ALC.Range("H77").Value = 3695.4119
$ws_ALC.Range("I77").Value = 3617.8
$ws_ALC.Range("J77").Value = 3806.2856
$ws_ALC.Range("K77").Value = 18089
$ws_ALC.Range("L77").Value = 19031.428
$ws_ALC.Range("M77").Value = -13409
$ws_ALC.Range("N77").Value = -28391.428
$ws_ALC.Range("H138").Value = 10875470
$ws_ALC.Range("I138").Value = 5621.5557
$ws_ALC.Range("J138").Value = 13519487
$ws_ALC.Range("K138").Value = 16864.6671
$ws_ALC.Range("L138").Value = 40558461
$ws_ALC.Range("M138").Value = -11724.6671
$ws_ALC.Range("N138").Value = -40568741

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H25").Value = 3471.6667
$ws_ARM.Range("J25").Value = 2499
$ws_ARM.Range("L25").Value = 2499
$ws_ARM.Range("N25").Value = -3303
$ws_ARM.Range("H61").Value = 2776.875
$ws_ARM.Range("I61").Value = 1725
$ws_ARM.Range("J61").Value = 5932.5
$ws_ARM.Range("K61").Value = 1725
$ws_ARM.Range("L61").Value = 5932.5
$ws_ARM.Range("M61").Value = -1513
$ws_ARM.Range("N61").Value = -6356.5
$ws_ARM.Range("H132").Value = 3415.8064
$ws_ARM.Range("I132").Value = 3009.647
$ws_ARM.Range("J132").Value = 3909
$ws_ARM.Range("K132").Value = 9028.940999999999
$ws_ARM.Range("L132").Value = 11727
$ws_ARM.Range("M132").Value = -6498.940999999999
$ws_ARM.Range("N132").Value = -16787
$ws_ARM.Range("H136").Value = 2776.875
$ws_ARM.Range("I136").Value = 1725
$ws_ARM.Range("J136").Value = 5932.5
$ws_ARM.Range("K136").Value = 5175
$ws_ARM.Range("L136").Value = 17797.5
$ws_ARM.Range("M136").Value = -2625
$ws_ARM.Range("N136").Value = -22897.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 5775.273
$ws_BSM.Range("I86").Value = 6763.7334
$ws_BSM.Range("J86").Value = 3657.1428
$ws_BSM.Range("K86").Value = 6763.7334
$ws_BSM.Range("L86").Value = 3657.1428
$ws_BSM.Range("M86").Value = -5640.7334
$ws_BSM.Range("N86").Value = -5903.1428
$ws_BSM.Range("H89").Value = 5775.273
$ws_BSM.Range("I89").Value = 6763.7334
$ws_BSM.Range("J89").Value = 3657.1428
$ws_BSM.Range("K89").Value = 33818.667
$ws_BSM.Range("L89").Value = 18285.714
$ws_BSM.Range("M89").Value = -28202.667
$ws_BSM.Range("N89").Value = -29517.714
$ws_BSM.Range("H105").Value = 3272.3809
$ws_BSM.Range("I105").Value = 2981.3333
$ws_BSM.Range("K105").Value = 2981.3333
$ws_BSM.Range("M105").Value = -1234.3333
$ws_BSM.Range("H134").Value = 3195
$ws_BSM.Range("I134").Value = 1804
$ws_BSM.Range("J134").Value = 5977
$ws_BSM.Range("K134").Value = 5412
$ws_BSM.Range("L134").Value = 17931
$ws_BSM.Range("M134").Value = -2877
$ws_BSM.Range("N134").Value = -23001

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H20").Value = 55599.4
$ws_CRP.Range("J20").Value = 55599.4
$ws_CRP.Range("L20").Value = 55599.4
$ws_CRP.Range("N20").Value = -56071.4
$ws_CRP.Range("H30").Value = 55599.4
$ws_CRP.Range("J30").Value = 55599.4
$ws_CRP.Range("L30").Value = 55599.4
$ws_CRP.Range("N30").Value = -55781.4
$ws_CRP.Range("H31").Value = 1856.35
$ws_CRP.Range("I31").Value = 1231.0769
$ws_CRP.Range("K31").Value = 1231.0769
$ws_CRP.Range("M31").Value = -936.0769
$ws_CRP.Range("H34").Value = 1856.35
$ws_CRP.Range("I34").Value = 1231.0769
$ws_CRP.Range("K34").Value = 1231.0769
$ws_CRP.Range("M34").Value = -1029.0769
$ws_CRP.Range("H58").Value = 794.2459
$ws_CRP.Range("I58").Value = 386.4634
$ws_CRP.Range("J58").Value = 1630.2
$ws_CRP.Range("K58").Value = 386.4634
$ws_CRP.Range("L58").Value = 1630.2
$ws_CRP.Range("M58").Value = -183.4634
$ws_CRP.Range("N58").Value = -2036.2
$ws_CRP.Range("H99").Value = 10420335
$ws_CRP.Range("I99").Value = 15630003
$ws_CRP.Range("J99").Value = 1000
$ws_CRP.Range("K99").Value = 15630003
$ws_CRP.Range("L99").Value = 1000
$ws_CRP.Range("M99").Value = -15628505
$ws_CRP.Range("N99").Value = -3996
$ws_CRP.Range("H126").Value = 10420335
$ws_CRP.Range("I126").Value = 15630003
$ws_CRP.Range("J126").Value = 1000
$ws_CRP.Range("K126").Value = 46890009
$ws_CRP.Range("L126").Value = 3000
$ws_CRP.Range("M126").Value = -46887539
$ws_CRP.Range("N126").Value = -7940
$ws_CRP.Range("H128").Value = 55599.4
$ws_CRP.Range("J128").Value = 55599.4
$ws_CRP.Range("L128").Value = 55599.4
$ws_CRP.Range("N128").Value = -65559.39999999999
$ws_CRP.Range("H136").Value = 794.2459
$ws_CRP.Range("I136").Value = 386.4634
$ws_CRP.Range("J136").Value = 1630.2
$ws_CRP.Range("K136").Value = 1159.3902
$ws_CRP.Range("L136").Value = 4890.6
$ws_CRP.Range("M136").Value = 1390.6098
$ws_CRP.Range("N136").Value = -9990.6

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 968.5278
$ws_CUL.Range("I5").Value = 431.375
$ws_CUL.Range("J5").Value = 1398.25
$ws_CUL.Range("K5").Value = 1294.125
$ws_CUL.Range("L5").Value = 4194.75
$ws_CUL.Range("M5").Value = -1182.125
$ws_CUL.Range("N5").Value = -4418.75
$ws_CUL.Range("H131").Value = 2382.1428
$ws_CUL.Range("I131").Value = 391.66666
$ws_CUL.Range("K131").Value = 1174.99998
$ws_CUL.Range("M131").Value = 3865.00002
$ws_CUL.Range("H132").Value = 832.26666
$ws_CUL.Range("I132").Value = 594.1429000000001
$ws_CUL.Range("J132").Value = 1040.625
$ws_CUL.Range("K132").Value = 5347.2861
$ws_CUL.Range("L132").Value = 9365.625
$ws_CUL.Range("M132").Value = -2817.2861
$ws_CUL.Range("N132").Value = -14425.625
$ws_CUL.Range("H133").Value = 10000
$ws_CUL.Range("I133").Value = 9000
$ws_CUL.Range("J133").Value = 13000
$ws_CUL.Range("K133").Value = 27000
$ws_CUL.Range("L133").Value = 39000
$ws_CUL.Range("M133").Value = -21940
$ws_CUL.Range("N133").Value = -49120
$ws_CUL.Range("H135").Value = 968.5278
$ws_CUL.Range("I135").Value = 431.375
$ws_CUL.Range("J135").Value = 1398.25
$ws_CUL.Range("K135").Value = 3882.375
$ws_CUL.Range("L135").Value = 12584.25
$ws_CUL.Range("M135").Value = -1347.375
$ws_CUL.Range("N135").Value = -17654.25
$ws_CUL.Range("H136").Value = 4015.862
$ws_CUL.Range("I136").Value = 948.2308
$ws_CUL.Range("J136").Value = 6508.3125
$ws_CUL.Range("K136").Value = 2844.6924
$ws_CUL.Range("L136").Value = 19524.9375
$ws_CUL.Range("M136").Value = 2255.3076
$ws_CUL.Range("N136").Value = -29724.9375
$ws_CUL.Range("H137").Value = 64555
$ws_CUL.Range("I137").Value = 4372.4165
$ws_CUL.Range("J137").Value = 208993.2
$ws_CUL.Range("K137").Value = 13117.2495
$ws_CUL.Range("L137").Value = 626979.6000000001
$ws_CUL.Range("M137").Value = -8017.249500000002
$ws_CUL.Range("N137").Value = -637179.6000000001
$ws_CUL.Range("H138").Value = 3317.9583
$ws_CUL.Range("I138").Value = 740.06665
$ws_CUL.Range("J138").Value = 7614.4443
$ws_CUL.Range("K138").Value = 2220.19995
$ws_CUL.Range("L138").Value = 22843.3329
$ws_CUL.Range("M138").Value = 2919.80005
$ws_CUL.Range("N138").Value = -33123.3329

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 2279
$ws_GSM.Range("I80").Value = 2333.3333
$ws_GSM.Range("J80").Value = 2069.4285
$ws_GSM.Range("K80").Value = 2333.3333
$ws_GSM.Range("L80").Value = 2069.4285
$ws_GSM.Range("M80").Value = -1335.3333
$ws_GSM.Range("N80").Value = -4065.4285
$ws_GSM.Range("H83").Value = 2279
$ws_GSM.Range("I83").Value = 2333.3333
$ws_GSM.Range("J83").Value = 2069.4285
$ws_GSM.Range("K83").Value = 11666.6665
$ws_GSM.Range("L83").Value = 10347.1425
$ws_GSM.Range("M83").Value = -6674.666499999999
$ws_GSM.Range("N83").Value = -20331.1425

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 3450
$ws_LTW.Range("I7").Value = 2500
$ws_LTW.Range("J7").Value = 3555.5557
$ws_LTW.Range("K7").Value = 2500
$ws_LTW.Range("L7").Value = 3555.5557
$ws_LTW.Range("M7").Value = -2388
$ws_LTW.Range("N7").Value = -3779.5557
$ws_LTW.Range("H40").Value = 4686
$ws_LTW.Range("I40").Value = 3395
$ws_LTW.Range("J40").Value = 4884.615
$ws_LTW.Range("K40").Value = 3395
$ws_LTW.Range("L40").Value = 4884.615
$ws_LTW.Range("M40").Value = -3259
$ws_LTW.Range("N40").Value = -5156.615
$ws_LTW.Range("H46").Value = 1009.25
$ws_LTW.Range("I46").Value = 575
$ws_LTW.Range("J46").Value = 1269.8
$ws_LTW.Range("K46").Value = 575
$ws_LTW.Range("L46").Value = 1269.8
$ws_LTW.Range("M46").Value = -387
$ws_LTW.Range("N46").Value = -1645.8
$ws_LTW.Range("H61").Value = 1471.2858
$ws_LTW.Range("I61").Value = 1674.75
$ws_LTW.Range("J61").Value = 1200
$ws_LTW.Range("K61").Value = 1674.75
$ws_LTW.Range("L61").Value = 1200
$ws_LTW.Range("M61").Value = -1472.75
$ws_LTW.Range("N61").Value = -1604
$ws_LTW.Range("H113").Value = 1471.2858
$ws_LTW.Range("I113").Value = 1674.75
$ws_LTW.Range("J113").Value = 1200
$ws_LTW.Range("K113").Value = 1674.75
$ws_LTW.Range("L113").Value = 1200
$ws_LTW.Range("M113").Value = 495.25
$ws_LTW.Range("N113").Value = -5540
$ws_LTW.Range("H122").Value = 3625.3684
$ws_LTW.Range("I122").Value = 2220.5
$ws_LTW.Range("J122").Value = 4000
$ws_LTW.Range("K122").Value = 6661.5
$ws_LTW.Range("L122").Value = 12000
$ws_LTW.Range("M122").Value = -4211.5
$ws_LTW.Range("N122").Value = -16900
$ws_LTW.Range("H126").Value = 3450
$ws_LTW.Range("I126").Value = 2500
$ws_LTW.Range("J126").Value = 3555.5557
$ws_LTW.Range("K126").Value = 7500
$ws_LTW.Range("L126").Value = 10666.6671
$ws_LTW.Range("M126").Value = -5030
$ws_LTW.Range("N126").Value = -15606.6671
$ws_LTW.Range("H128").Value = 67900
$ws_LTW.Range("J128").Value = 67900
$ws_LTW.Range("L128").Value = 67900
$ws_LTW.Range("N128").Value = -77860
$ws_LTW.Range("H130").Value = 0
$ws_LTW.Range("J130").Value = 0
$ws_LTW.Range("L130").Value = 0
$ws_LTW.Range("N130").ClearContents()
$ws_LTW.Range("H131").Value = 33150
$ws_LTW.Range("J131").Value = 0
$ws_LTW.Range("L131").Value = 0
$ws_LTW.Range("N131").ClearContents()
$ws_LTW.Range("H132").Value = 4522.1763
$ws_LTW.Range("I132").Value = 3967.4443
$ws_LTW.Range("J132").Value = 5146.25
$ws_LTW.Range("K132").Value = 11902.3329
$ws_LTW.Range("L132").Value = 15438.75
$ws_LTW.Range("M132").Value = -9372.332900000001
$ws_LTW.Range("N132").Value = -20498.75

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H126").Value = 125812.625
$ws_WVR.Range("I126").Value = 167133.5
$ws_WVR.Range("J126").Value = 1850
$ws_WVR.Range("K126").Value = 501400.5
$ws_WVR.Range("L126").Value = 5550
$ws_WVR.Range("M126").Value = -498930.5
$ws_WVR.Range("N126").Value = -10490
$ws_WVR.Range("H133").Value = 79215
$ws_WVR.Range("J133").Value = 79215
$ws_WVR.Range("L133").Value = 79215
$ws_WVR.Range("N133").Value = -89335
